$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New results ("2021/05/30 14:50") for the two updated data points.
$ws.Range("B17").Value = 0.3952
$ws.Range("B18").Value = 0.444

# B32 holds =AVERAGE(B2:B31); it recalculates automatically from the
# updated inputs above, moving from 0.38801333333333327 to 0.38927999999999996.

# Reposition the view the way it was left after entering the new values:
# scrolled down so row 4 is the top visible row, with B18 the active cell.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B18").Select()
